# Trade #61 closed at 2026-02-17 15:44:19 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status roll-up numbers for the
# MarketMaking strategy and appends the newly closed trade (#61) as
# row 62 on both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet - roll-up metrics
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.31   # Current Capital
$summary.Range("B4").Value = 0.31      # Total P&L $
$summary.Range("B5").Value = 0.1       # Total P&L %
$summary.Range("B6").Value = 61        # Total Trades
$summary.Range("B8").Value = 33        # Losing Trades
$summary.Range("B9").Value = 31.15     # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.31     # Capital
$status.Range("D4").Value = 61         # Trades
$status.Range("E4").Value = 0.31       # P&L $
$status.Range("F4").Value = 0.31       # P&L %
$status.Range("G4").Value = 31.15      # Win Rate %

# ---------------------------------------------------------------------
# Helper: append the new trade record (#61) as row 62 on a trades sheet.
# Date/time columns must stay as plain text (not auto-converted to Excel
# date/time serials), so we stage them with a Text number format and
# then clear the format again to avoid leaving a stray style behind.
# ---------------------------------------------------------------------
function Add-Trade61Row($ws) {
    $row = 62

    $ws.Cells.Item($row, 1).Value = 61

    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 2).ClearFormats()

    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = "15:44:12"
    $ws.Cells.Item($row, 3).ClearFormats()

    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.86
    $ws.Cells.Item($row, 7).Value = 0.83
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -3.4884
    $ws.Cells.Item($row, 10).Value = -0.03
    $ws.Cells.Item($row, 11).Value = 100.31
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.11
}

# ---------------------------------------------------------------------
# 3. All Trades sheet - append row 62
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade61Row $allTrades

# ---------------------------------------------------------------------
# 4. MarketMaking sheet - append the same row 62
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade61Row $marketMaking
